$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Column I width: 14 -> 10 ---
$ws.Columns.Item(9).ColumnWidth = 10

# --- 2. Class statistics scalar updates (rows 6,7,9) ---
$ws.Range("L6").Value = 90
$ws.Range("L7").Value = 0
$ws.Range("L9").Value = "'40.5%"

# --- 3. Swap "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System" ---
$gSwapRows = @(2,3,20,21,22,24,39,40,41,43,58,59,60,62,77,78,95,96,113,114,131,132,149,150,167,168,169,171,186,187,188,190,205,206,207,209)
foreach ($r in $gSwapRows) {
    $ws.Range("G" + $r).Value = "dnasr281@gmail.com, System"
}

# --- 4. Group statistics updates for rows 15-26 (O,P,R,S) ---
$ws.Range("O15").Value = 7
$ws.Range("P15").Value = 0
$ws.Range("R15").Value = "'38.9%"
$ws.Range("S15").Value = "'76.7%"
$ws.Range("O16").Value = 8
$ws.Range("P16").Value = 0
$ws.Range("R16").Value = "'42.1%"
$ws.Range("S16").Value = "'71.8%"
$ws.Range("O17").Value = 8
$ws.Range("P17").Value = 0
$ws.Range("R17").Value = "'42.1%"
$ws.Range("S17").Value = "'59.9%"
$ws.Range("O18").Value = 8
$ws.Range("P18").Value = 0
$ws.Range("R18").Value = "'42.1%"
$ws.Range("S18").Value = "'83.9%"
$ws.Range("O19").Value = 7
$ws.Range("P19").Value = 0
$ws.Range("R19").Value = "'38.9%"
$ws.Range("S19").Value = "'88.0%"
$ws.Range("O20").Value = 7
$ws.Range("P20").Value = 0
$ws.Range("R20").Value = "'38.9%"
$ws.Range("S20").Value = "'90.8%"
$ws.Range("O21").Value = 7
$ws.Range("P21").Value = 0
$ws.Range("R21").Value = "'38.9%"
$ws.Range("S21").Value = "'89.7%"
$ws.Range("O22").Value = 7
$ws.Range("P22").Value = 0
$ws.Range("R22").Value = "'38.9%"
$ws.Range("S22").Value = "'90.5%"
$ws.Range("O23").Value = 7
$ws.Range("P23").Value = 0
$ws.Range("R23").Value = "'38.9%"
$ws.Range("S23").Value = "'67.6%"
$ws.Range("O24").Value = 8
$ws.Range("P24").Value = 0
$ws.Range("R24").Value = "'42.1%"
$ws.Range("S24").Value = "'70.4%"
$ws.Range("O25").Value = 8
$ws.Range("P25").Value = 0
$ws.Range("R25").Value = "'42.1%"
$ws.Range("S25").Value = "'74.6%"
$ws.Range("O26").Value = 8
$ws.Range("P26").Value = 0
$ws.Range("R26").Value = "'42.1%"
$ws.Range("S26").Value = "'71.1%"

# --- 5. Newly recorded sessions (style pink->green, content updates) ---
$recorded = @{
    16 = "21/27"
    35 = "18/31"
    54 = "13/19"
    73 = "18/21"
    91 = "25/31"
    109 = "25/28"
    127 = "21/29"
    145 = "28/33"
    163 = "24/30"
    182 = "20/27"
    201 = "23/29"
    220 = "24/29"
}
foreach ($r in $recorded.Keys) {
    $src = $ws.Range("A2:I2")
    $src.Copy()
    $dst = $ws.Range("A" + $r + ":I" + $r)
    $dst.PasteSpecial(-4122)
    $ws.Range("G" + $r).Value = "dnasr281@gmail.com"
    $ws.Range("H" + $r).Value = $recorded[$r]
    $ws.Range("I" + $r).Value = "Recorded"
}
$excel.CutCopyMode = $false

Write-Output "edit complete"
